$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 76 (shifts existing rows 76..166 down to 77..167)
$ws.Rows.Item(76).Insert()

# Populate the newly inserted row 76 with the new weekly price record
$row = 76
$ws.Cells.Item($row, 1).Value2  = 11
$ws.Cells.Item($row, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value2  = "Bíobío"
$ws.Cells.Item($row, 4).Value2  = 44587
$ws.Cells.Item($row, 5).Value2  = 8
$ws.Cells.Item($row, 6).Value2  = 100112045
$ws.Cells.Item($row, 7).Value2  = "Zapallo"
$ws.Cells.Item($row, 8).Value2  = "Camote"
$ws.Cells.Item($row, 9).Value2  = "1a (cosecha)"
$ws.Cells.Item($row, 10).Value2 = 450
$ws.Cells.Item($row, 11).Value2 = 500
$ws.Cells.Item($row, 12).Value2 = 550
$ws.Cells.Item($row, 13).Value2 = 528
$ws.Cells.Item($row, 14).Value2 = "`$/kilo"
$ws.Cells.Item($row, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item($row, 16).Value2 = 528
$ws.Cells.Item($row, 17).Value2 = 1
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
